$d = $word.ActiveDocument

# The last paragraph ends with the hidden "_GoBack" bookmark markers
# (Word's "last edit location" bookmark) sitting right after the
# existing text, before the paragraph mark. A plain text insertion at
# the end of the story lands *before* those markers (merging into the
# existing last run), while a structural paragraph-mark insertion at
# the end of the story lands *after* them. To reproduce the target
# edit - a brand-new numbered-list paragraph "Merge " that itself ends
# up owning the trailing bookmark markers - we:
#   1) append the new paragraph's text onto the end of the story first
#      (this naturally lands before the bookmark markers), then
#   2) split it off into its own paragraph by inserting a paragraph
#      mark immediately before the text we just added.
$end = $d.Content
$end.Collapse(0)
$end.InsertAfter("Merge ")

$split = $d.Content
[void]$split.Find.Execute("Merge ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$split.Collapse(1)
$split.InsertParagraphBefore()
